# Tidy command in pptx
#
# 1) Date placeholder ("datetimeFigureOut" field) cached text changes
#    from 2023/1/12 -> 2023/1/28 on the slide master and every slide
#    layout.
# 2) Slide 2 ("git 對象庫" content slide), bullet "5." text tidy-up:
#    "查询目前tree里面的文件" -> "查询暂存的文件"

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1) Update the cached date field text wherever it appears
#    (slide master + all custom layouts) by locating the
#    placeholder of type ppPlaceholderDate (16).
# ---------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        $phType = -1
        try { $phType = $shape.PlaceholderFormat.Type } catch {}
        if ($phType -eq 16) {
            if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
                $tr = $shape.TextFrame.TextRange
                if ($tr.Text -eq "2023/1/12") {
                    $tr.Text = "2023/1/28"
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------
# 2) Slide 2, paragraph "5. ...": retype the run(s) that read
#    "查询目前tree里面的" into "查询暂存的", leaving the trailing
#    "文件" run untouched.
# ---------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$shape = $slide2.Shapes.Item(1)
$textRange = $shape.TextFrame.TextRange

for ($i = 1; $i -le $textRange.Paragraphs().Count; $i++) {
    $para = $textRange.Paragraphs($i)
    if ($para.Text.StartsWith("5. ")) {
        # "5. " (3 chars) is left alone; the next 11 characters
        # ("查询目前tree里面的") are retyped, the final "文件" stays.
        $target = $para.Characters(4, 11)
        $target.Text = "查询暂存的"
        break
    }
}
